# all resitors & capacitors from 0603 to 0402, new eurocircuits basket B2206781
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the header time string (B1-ish info cell E24 holds "17:42" -> "11:49")
$ws.Range("E24").Value = "11:49"

# --- Capacitors: rows 2-5 (C1,C4,C5 / C2,C3 / C6,C7 / C8,C10) ---
# Footprint column D: 0603_CAP_SMALL -> 0402_CAP
$ws.Range("D2").Value = "0402_CAP"
$ws.Range("D3").Value = "0402_CAP"
$ws.Range("D4").Value = "0402_CAP"
$ws.Range("D5").Value = "0402_CAP"

# Manufacturer column E
$ws.Range("E2").Value = "Wurth Electronics"
$ws.Range("E3").Value = "Walsin Technologies"
$ws.Range("E4").Value = "Wurth Electronics"
$ws.Range("E5").Value = "Taiyo Yuden"

# Manufacturer Part Number column F
$ws.Range("F2").Value = "885012105016"
$ws.Range("F3").Value = "0402N180F500CT"
$ws.Range("F4").Value = "885012005058"
$ws.Range("F5").Value = "JMK105BJ105KP-F"

# --- Resistors: rows 12-18 (R1,R3,R4,R9,R14 / R2 / R5,R7 / R6,R8 / R10 / R12 / R15) ---
# Footprint column D: 0603_res_SMALL -> 0402_res (or 0402_RES for rows 14 & 18)
$ws.Range("D12").Value = "0402_res"
$ws.Range("D13").Value = "0402_res"
$ws.Range("D14").Value = "0402_RES"
$ws.Range("D15").Value = "0402_res"
$ws.Range("D16").Value = "0402_res"
$ws.Range("D17").Value = "0402_res"
$ws.Range("D18").Value = "0402_RES"

# Manufacturer column E
$ws.Range("E12").Value = "Vishay"
$ws.Range("E13").Value = "Vishay"
$ws.Range("E14").Value = "Vishay Semiconductors"
$ws.Range("E15").Value = "Vishay Semiconductors"
$ws.Range("E16").Value = "Vishay"
$ws.Range("E17").Value = "Vishay"
$ws.Range("E18").Value = "Vishay"

# Manufacturer Part Number column F
$ws.Range("F12").Value = "CRCW040210K0FKEDC"
$ws.Range("F13").Value = "CRCW0402680RFKEDC"
$ws.Range("F14").Value = "CRCW0402470RFKEDC"
$ws.Range("F15").Value = "CRCW040233R0FKEDC"
$ws.Range("F16").Value = "CRCW04021K50FKEDC"
$ws.Range("F17").Value = "CRCW04020000Z0EDC"
$ws.Range("F18").Value = "CRCW040218K0FKED"

# Row-height autofit follow-on: longer wrapped manufacturer names (e.g.
# "Walsin Technologies", "Vishay Semiconductors") push these rows to wrap
# onto a second line in the narrow Manufacturer column, so Excel grows the
# row height from 15 to 21 points, same as rows 9/10/20/21 already do.
$ws.Rows("3").RowHeight = 21
$ws.Rows("14").RowHeight = 21
$ws.Rows("15").RowHeight = 21

$wb.Save()
